$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 374.8
$ws.Range("I33").Value = 267.57144
$ws.Range("K33").Value = 267.57144
$ws.Range("M33").Value = -38.57144
$ws.Range("H41").Value = 1752.7742
$ws.Range("I41").Value = 1512.2142
$ws.Range("J41").Value = 1950.8823
$ws.Range("K41").Value = 1512.2142
$ws.Range("L41").Value = 1950.8823
$ws.Range("M41").Value = -1072.2142
$ws.Range("N41").Value = -2830.8823
$ws.Range("H108").Value = 42796.09
$ws.Range("J108").Value = 42796.09
$ws.Range("L108").Value = 42796.09
$ws.Range("N108").Value = -50476.09
$ws.Range("H110").Value = 63970
$ws.Range("J110").Value = 63970
$ws.Range("L110").Value = 63970
$ws.Range("N110").Value = -72150
$ws.Range("H123").Value = 93863.336
$ws.Range("J123").Value = 93863.336
$ws.Range("L123").Value = 93863.336
$ws.Range("N123").Value = -103663.336
$ws.Range("H125").Value = 832.9375
$ws.Range("I125").Value = 532.6667
$ws.Range("K125").Value = 4794.0003
$ws.Range("M125").Value = -2334.0003
$ws.Range("H132").Value = 1990.1875
$ws.Range("J132").Value = 13002
$ws.Range("L132").Value = 39006
$ws.Range("N132").Value = -44066
$ws.Range("H133").Value = 93945.45
$ws.Range("J133").Value = 93945.45
$ws.Range("L133").Value = 93945.45
$ws.Range("N133").Value = -104065.45
$ws.Range("H138").Value = 2238.232
$ws.Range("I138").Value = 1599.5333
$ws.Range("K138").Value = 4798.5999
$ws.Range("M138").Value = 341.4000999999998
$ws.Range("H139").Value = 99988.336
$ws.Range("J139").Value = 99988.336
$ws.Range("L139").Value = 99988.336
$ws.Range("N139").Value = -110268.336
$ws.Range("H140").Value = 91980
$ws.Range("J140").Value = 91980
$ws.Range("L140").Value = 91980
$ws.Range("N140").Value = -102340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 512500
$ws.Range("I13").Value = 512500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 512500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -512356
$ws.Range("H52").Value = 68411.125
$ws.Range("J52").Value = 68411.125
$ws.Range("L52").Value = 68411.125
$ws.Range("N52").Value = -69047.125
$ws.Range("H74").Value = 71430.13
$ws.Range("I74").Value = 113680.445
$ws.Range("K74").Value = 113680.445
$ws.Range("M74").Value = -112806.445
$ws.Range("H77").Value = 71430.13
$ws.Range("I77").Value = 113680.445
$ws.Range("K77").Value = 568402.2250000001
$ws.Range("M77").Value = -564034.2250000001
$ws.Range("H107").Value = 49132.145
$ws.Range("J107").Value = 49132.145
$ws.Range("L107").Value = 49132.145
$ws.Range("N107").Value = -56812.145
$ws.Range("H108").Value = 89992.5
$ws.Range("J108").Value = 89992.5
$ws.Range("L108").Value = 89992.5
$ws.Range("N108").Value = -97672.5
$ws.Range("H118").Value = 89990
$ws.Range("J118").Value = 89990
$ws.Range("L118").Value = 89990
$ws.Range("N118").Value = -93304
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 64940
$ws.Range("J6").Value = 64940
$ws.Range("L6").Value = 64940
$ws.Range("N6").Value = -65166
$ws.Range("H55").Value = 38757.25
$ws.Range("J55").Value = 38773.332
$ws.Range("L55").Value = 38773.332
$ws.Range("N55").Value = -39319.332
$ws.Range("H117").Value = 83737.5
$ws.Range("J117").Value = 83737.5
$ws.Range("L117").Value = 83737.5
$ws.Range("N117").Value = -92915.5
$ws.Range("H119").Value = 85866.75
$ws.Range("J119").Value = 85866.75
$ws.Range("L119").Value = 85866.75
$ws.Range("N119").Value = -95542.75
$ws.Range("H122").Value = 85328.336
$ws.Range("J122").Value = 85328.336
$ws.Range("L122").Value = 85328.336
$ws.Range("N122").Value = -95128.336
$ws.Range("H132").Value = 35221.777
$ws.Range("J132").Value = 35221.777
$ws.Range("L132").Value = 35221.777
$ws.Range("N132").Value = -45341.777
$ws.Range("H134").Value = 2092.7368
$ws.Range("I134").Value = 1497.2307
$ws.Range("J134").Value = 3383
$ws.Range("K134").Value = 4491.6921
$ws.Range("L134").Value = 10149
$ws.Range("M134").Value = -1956.6921
$ws.Range("N134").Value = -15219
$ws.Range("H135").Value = 101278.57
$ws.Range("J135").Value = 101278.57
$ws.Range("L135").Value = 101278.57
$ws.Range("N135").Value = -111418.57
$ws.Range("H138").Value = 95991.664
$ws.Range("J138").Value = 95991.664
$ws.Range("L138").Value = 95991.664
$ws.Range("N138").Value = -106271.664
$ws.Range("H140").Value = 70106.125
$ws.Range("J140").Value = 70106.125
$ws.Range("L140").Value = 70106.125
$ws.Range("N140").Value = -80466.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1737.625
$ws.Range("I58").Value = 1460.1818
$ws.Range("J58").Value = 2348
$ws.Range("K58").Value = 1460.1818
$ws.Range("L58").Value = 2348
$ws.Range("M58").Value = -1257.1818
$ws.Range("N58").Value = -2754
$ws.Range("H107").Value = 747.58826
$ws.Range("I107").Value = 775.44446
$ws.Range("K107").Value = 775.44446
$ws.Range("M107").Value = 1144.55554
$ws.Range("H108").Value = 42642.715
$ws.Range("J108").Value = 42642.715
$ws.Range("L108").Value = 42642.715
$ws.Range("N108").Value = -50322.715
$ws.Range("H114").Value = 47003
$ws.Range("J114").Value = 47003
$ws.Range("L114").Value = 47003
$ws.Range("N114").Value = -55681
$ws.Range("H116").Value = 43103.5
$ws.Range("J116").Value = 43103.5
$ws.Range("L116").Value = 43103.5
$ws.Range("N116").Value = -52281.5
$ws.Range("H118").Value = 83874.89
$ws.Range("J118").Value = 83874.89
$ws.Range("L118").Value = 83874.89
$ws.Range("N118").Value = -87188.89
$ws.Range("H119").Value = 98491.664
$ws.Range("J119").Value = 98491.664
$ws.Range("L119").Value = 98491.664
$ws.Range("N119").Value = -108167.664
$ws.Range("H132").Value = 1978.909
$ws.Range("I132").Value = 1508.0588
$ws.Range("K132").Value = 4524.1764
$ws.Range("M132").Value = -1994.1764
$ws.Range("H136").Value = 1737.625
$ws.Range("I136").Value = 1460.1818
$ws.Range("J136").Value = 2348
$ws.Range("K136").Value = 4380.5454
$ws.Range("L136").Value = 7044
$ws.Range("M136").Value = -1830.5454
$ws.Range("N136").Value = -12144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2712.5
$ws.Range("I70").Value = 425
$ws.Range("K70").Value = 1275
$ws.Range("M70").Value = -960
$ws.Range("H73").Value = 2712.5
$ws.Range("I73").Value = 425
$ws.Range("K73").Value = 1275
$ws.Range("M73").Value = -183

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 64943.332
$ws.Range("J108").Value = 64943.332
$ws.Range("L108").Value = 64943.332
$ws.Range("N108").Value = -72623.33199999999
$ws.Range("H119").Value = 69407.28999999999
$ws.Range("J119").Value = 69526.664
$ws.Range("L119").Value = 69526.664
$ws.Range("N119").Value = -79202.664
$ws.Range("H135").Value = 69868
$ws.Range("J135").Value = 69868
$ws.Range("L135").Value = 69868
$ws.Range("N135").Value = -80008
$ws.Range("H141").Value = 119775.29
$ws.Range("J141").Value = 134606.75
$ws.Range("L141").Value = 134606.75
$ws.Range("N141").Value = -144966.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 39660.8
$ws.Range("J121").Value = 39663.75
$ws.Range("L121").Value = 39663.75
$ws.Range("N121").Value = -43157.75
$ws.Range("H123").Value = 84994.44500000001
$ws.Range("J123").Value = 84994.44500000001
$ws.Range("L123").Value = 84994.44500000001
$ws.Range("N123").Value = -94794.44500000001
$ws.Range("H129").Value = 100000
$ws.Range("J129").Value = 100000
$ws.Range("L129").Value = 100000
$ws.Range("N129").Value = -110000
$ws.Range("H136").Value = 5085.6113
$ws.Range("I136").Value = 4643.421
$ws.Range("K136").Value = 13930.263
$ws.Range("M136").Value = -11380.263
$ws.Range("H138").Value = 105943
$ws.Range("J138").Value = 121590.664
$ws.Range("L138").Value = 121590.664
$ws.Range("N138").Value = -131870.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 89897
$ws.Range("J121").Value = 89897
$ws.Range("L121").Value = 89897
$ws.Range("N121").Value = -93391
$ws.Range("H127").Value = 61496.5
$ws.Range("J127").Value = 61496.5
$ws.Range("L127").Value = 61496.5
$ws.Range("N127").Value = -71416.5
$ws.Range("H136").Value = 2361.3333
$ws.Range("I136").Value = 1684.409
$ws.Range("K136").Value = 5053.227000000001
$ws.Range("M136").Value = -2503.227000000001
